$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price (D) / volume-change (E) updates ---
$ws.Range("D2").Value = '71.046.62'
$ws.Range("E2").Value = '  -2.65%  '

$ws.Range("D3").Value = '3.845.49'
$ws.Range("E3").Value = '  -3.64%  '

$ws.Range("E4").Value = '  +0.38%  '

$ws.Range("D5").Value = '599.31'
$ws.Range("E5").Value = '  -2.61%  '

$ws.Range("D6").Value = '182.64'
$ws.Range("E6").Value = '  +10.16%  '

$ws.Range("D7").Value = '0.666'
$ws.Range("E7").Value = '  -2.97%  '

$ws.Range("D8").Value = "'" + '1.00'
$ws.Range("E8").Value = '  +0.43%  '

$ws.Range("D9").Value = '0.755'
$ws.Range("E9").Value = '  -0.47%  '

$ws.Range("D10").Value = '0.177'
$ws.Range("E10").Value = '  +5.63%  '

$ws.Range("D11").Value = '55.33'
$ws.Range("E11").Value = '  -4.84%  '

$ws.Range("D12").Value = '0.0000315'
$ws.Range("E12").Value = '  -0.67%  '

$ws.Range("D13").Value = "'" + '11.30'
$ws.Range("E13").Value = '  +0.86%  '

$ws.Range("D14").Value = '4.483.89'
$ws.Range("E14").Value = '  -3.01%  '

$ws.Range("D15").Value = '3.867.13'
$ws.Range("E15").Value = '  -3.49%  '

$ws.Range("D16").Value = '20.38'
$ws.Range("E16").Value = '  -1.59%  '

$ws.Range("D17").Value = '13.72'
$ws.Range("E17").Value = '  -4.21%  '

$ws.Range("D18").Value = "'" + '1.20'
$ws.Range("E18").Value = '  -6.03%  '

$ws.Range("E19").Value = '  -2.40%  '

$ws.Range("D20").Value = '70.967.09'
$ws.Range("E20").Value = '  -2.50%  '

$ws.Range("D21").Value = '430.96'
$ws.Range("E21").Value = '  -2.45%  '

$ws.Range("D22").Value = '4.77'
$ws.Range("E22").Value = '  -3.21%  '

$ws.Range("D23").Value = '93.39'
$ws.Range("E23").Value = '  -3.33%  '

$ws.Range("D24").Value = '3.21'
$ws.Range("E24").Value = '  -5.56%  '

$ws.Range("D25").Value = '13.59'
$ws.Range("E25").Value = '  -6.95%  '

$ws.Range("D26").Value = '11.41'
$ws.Range("E26").Value = '  +0.38%  '

$ws.Range("D27").Value = '3.95'
$ws.Range("E27").Value = '  -6.03%  '

$ws.Range("D28").Value = '5.97'
$ws.Range("E28").Value = '  +0.25%  '

$ws.Range("D29").Value = '10.28'
$ws.Range("E29").Value = '  -2.57%  '

$ws.Range("D30").Value = '8.73'
$ws.Range("E30").Value = '  +12.09%  '

$ws.Range("D31").Value = '34.76'
$ws.Range("E31").Value = '  -4.04%  '

$ws.Range("D32").Value = '13.49'
$ws.Range("E32").Value = '  -2.98%  '

$ws.Range("D33").Value = '47.22'
$ws.Range("E33").Value = '  -4.02%  '

$ws.Range("E34").Value = '  -5.14%  '

$ws.Range("D35").Value = '0.0₃0973'
$ws.Range("E35").Value = '  +7.23%  '

$ws.Range("D38").Value = '0.422'
$ws.Range("E38").Value = '  -3.30%  '

$ws.Range("D39").Value = '3.39'
$ws.Range("E39").Value = '  +15.68%  '

$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.02%  '

$ws.Range("D41").Value = '0.144'
$ws.Range("E41").Value = '  -2.85%  '

$ws.Range("E42").Value = '  +0.21%  '

$ws.Range("D43").Value = '3.25'
$ws.Range("E43").Value = '  -7.95%  '

$ws.Range("D44").Value = '0.0466'
$ws.Range("E44").Value = '  -4.62%  '

$ws.Range("D45").Value = '10.04'
$ws.Range("E45").Value = '  -7.45%  '

$ws.Range("D46").Value = '2.69'
$ws.Range("E46").Value = '  +1.21%  '

$ws.Range("D47").Value = '0.142'
$ws.Range("E47").Value = '  -4.79%  '

$ws.Range("D51").Value = '0.000274'
$ws.Range("E51").Value = '  +3.33%  '

# --- Rows whose coin entry moved to a different rank position ---
$ws.Range("B36").Value = 'Bittensor'
$ws.Range("C36").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D36").Value = '637.13'
$ws.Range("E36").Value = '  +0.60%  '

$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = "'" + '68.00'
$ws.Range("E37").Value = '  -5.81%  '

$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").Value = '2.83'
$ws.Range("E48").Value = '  -15.50%  '

$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = '3.31'
$ws.Range("E49").Value = '  -2.71%  '

$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.873.84'
$ws.Range("E50").Value = '  -2.02%  '

